$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data (A:E -> B:F)
$ws.Columns.Item(1).Insert()

# Header for the new ID column, matching header style of the other headers (s="1")
$ws.Cells.Item(1, 1).Value = "ID"
$ws.Cells.Item(1, 2).Copy()
$ws.Cells.Item(1, 1).PasteSpecial(-4122)

# Fill in the ID labels for each data row
$ids = @("Hb 19", "Hb 20", "S 17", "Hb 50", "Hb 28", "Hb 30", "Hb 23", "Hb 26", "Hb 27")
for ($i = 0; $i -lt $ids.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $ids[$i]
}
